{"js": "// Replace the contents of the 5x3 lattice-multiplication table: each cell's\n// text (5 lines joined by manual line breaks, i.e. \"\\v\") is replaced with a\n// newly generated problem. The table's shape (5 rows x 3 columns, one\n// paragraph/run per cell, sz=32) is unchanged -- only the <w:t> text moves.\n\nconst newValues = [\n  \"81 x 19\\v  1    9\\v  ----\\v8|    |\\v1|    |\",\n  \"39 x 30\\v  3    0\\v  ----\\v3|    |\\v9|    |\",\n  \"96 x 29\\v  2    9\\v  ----\\v9|    |\\v6|    |\",\n  \"71 x 90\\v  9    0\\v  ----\\v7|    |\\v1|    |\",\n  \"37 x 33\\v  3    3\\v  ----\\v3|    |\\v7|    |\",\n  \"83 x 79\\v  7    9\\v  ----\\v8|    |\\v3|    |\",\n  \"94 x 85\\v  8    5\\v  ----\\v9|    |\\v4|    |\",\n  \"84 x 80\\v  8    0\\v  ----\\v8|    |\\v4|    |\",\n  \"68 x 32\\v  3    2\\v  ----\\v6|    |\\v8|    |\",\n  \"98 x 34\\v  3    4\\v  ----\\v9|    |\\v8|    |\",\n  \"92 x 93\\v  9    3\\v  ----\\v9|    |\\v2|    |\",\n  \"53 x 48\\v  4    8\\v  ----\\v5|    |\\v3|    |\",\n  \"58 x 92\\v  9    2\\v  ----\\v5|    |\\v8|    |\",\n  \"37 x 52\\v  5    2\\v  ----\\v3|    |\\v7|    |\",\n  \"86 x 80\\v  8    0\\v  ----\\v8|    |\\v6|    |\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columnCount = 3;\nconst rowCount = table.rowCount;\n\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const idx = r * columnCount + c;\n    if (idx >= newValues.length) continue;\n\n    const cell = table.getCell(r, c);\n    cell.body.paragraphs.load(\"items\");\n    await context.sync();\n\n    const para = cell.body.paragraphs.items[0];\n    const range = para.getRange();\n    range.insertText(newValues[idx], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the contents of the 5x3 lattice-multiplication table: each cell's\n# text (5 lines joined by manual line breaks, i.e. a vertical-tab \"`v\") is\n# replaced with a newly generated problem. The table's shape (5 rows x 3\n# columns, one paragraph/run per cell, sz=32) is unchanged -- only the\n# <w:t> text moves.\n\n$newValues = @(\n    \"81 x 19`v  1    9`v  ----`v8|    |`v1|    |\",\n    \"39 x 30`v  3    0`v  ----`v3|    |`v9|    |\",\n    \"96 x 29`v  2    9`v  ----`v9|    |`v6|    |\",\n    \"71 x 90`v  9    0`v  ----`v7|    |`v1|    |\",\n    \"37 x 33`v  3    3`v  ----`v3|    |`v7|    |\",\n    \"83 x 79`v  7    9`v  ----`v8|    |`v3|    |\",\n    \"94 x 85`v  8    5`v  ----`v9|    |`v4|    |\",\n    \"84 x 80`v  8    0`v  ----`v8|    |`v4|    |\",\n    \"68 x 32`v  3    2`v  ----`v6|    |`v8|    |\",\n    \"98 x 34`v  3    4`v  ----`v9|    |`v8|    |\",\n    \"92 x 93`v  9    3`v  ----`v9|    |`v2|    |\",\n    \"53 x 48`v  4    8`v  ----`v5|    |`v3|    |\",\n    \"58 x 92`v  9    2`v  ----`v5|    |`v8|    |\",\n    \"37 x 52`v  5    2`v  ----`v3|    |`v7|    |\",\n    \"86 x 80`v  8    0`v  ----`v8|    |`v6|    |\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$columnCount = 3\n$rowCount = $t.Rows.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $columnCount; $c++) {\n        $idx = (($r - 1) * $columnCount) + ($c - 1)\n        if ($idx -ge $newValues.Length) { continue }\n\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n    }\n}\n"}
